$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.253448963165283
$ws.Range("B1").Value = 2.267096519470215
$ws.Range("C1").Value = 4.50060510635376
$ws.Range("D1").Value = 2.90143346786499
$ws.Range("E1").Value = 1.370692491531372
